$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target data for the "Cluster name" / "Active cases" table (rows 2..35).
# This reflects: two new outbreak clusters added ("4518 Regis Aged Care Fawkner
# Outbreak" and "AW Window Transport Group Depot North Geelong Outbreak"),
# one new cluster added ("Community Kids Pascoe Vale Early Education Centre
# Pascoe Vale Outbreak"), one cluster removed ("Social Gathering Warrnambool
# 28 Sep Outbreak"), "Australian Lamb Colac East" replaced, and updated
# active-case counts throughout.
$data = @(
    @('139 Highett St Apartment Complex Richmond', 10),
    @('3175 The Bays Aged Care Facility Hastings', 17),
    @('3600 Belvedere Age Care Noble Park Outbreak', 21),
    @('3612 BlueCross Glengowrie Outbreak', 33),
    @('3684 Homestyle Aged Care Langford Grange Cranbourne East Outbreak', 23),
    @('4075 Ferndale Gardens Aged Care Services Bayswater North Outbreak', 17),
    @('4518 Regis Aged Care Fawkner Outbreak', 13),
    @('AW Window Transport Group Depot North Geelong Outbreak', 10),
    @('Bread Solutions Braeside Outbreak', 18),
    @('CS Square Caroline Springs Outbreak', 14),
    @('Cedar Meats Australia Brooklyn Outbreak', 10),
    @('Child''s Play Early Learning Centre Tarneit', 11),
    @('Community Kids Pascoe Vale Early Education Centre Pascoe Vale Outbreak', 20),
    @('Embracia Aged Care Reservoir Outbreak', 19),
    @('Guardian Childcare Caulfield Outbreak', 20),
    @('Hello Fresh Warehouse Ravenhall', 13),
    @('Inghams Enterprise Somerville Outbreak', 17),
    @('Kool Kidz Childcare Narre Warren', 10),
    @('Lantmannen Unibake Australia Mordialloc', 25),
    @('Launch Housing City Edge Crisis Accommodation South Melbourne', 11),
    @('Nido Early School Ascot Vale', 11),
    @('Northern Health Northern Hospital Epping Emergency Department Tier 1B', 41),
    @('Northern Health The Northern Hospital Epping', 14),
    @('Oceania Meat Processors Laverton North Outbreak', 16),
    @('Pick It Up Fitness Mulgrave Outbreak', 11),
    @('Robin Hood Inn Drouin West Outbreak', 42),
    @('St Vincents Hospital Emergency Department Melbourne', 41),
    @('Target Distribution Centre Truganina Outbreak', 19),
    @('The Royal Children''s Hospital Melbourne Emergency Department Parkville Tier 1B', 12),
    @('The Toolshed Bar Private Event Noojee', 17),
    @('Turosi Breakwater', 10),
    @('Visy Recycling Springvale', 21),
    @('Werribee Mercy Hospital Emergency Department', 25),
    @('Western Health Sunshine Hospital Emergency Department', 21)
)

# Keep the header row untouched, then (re)write the full data block, which
# both corrects existing rows and appends the new row 35 as needed.
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
